$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# Remove the "password" column (column E) entirely, shifting later columns left.
$ws.Columns("E").Delete()

# Update the active selection to match the post-edit state.
$ws.Range("I11").Select()
